$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 - Robert, Criterion 1 Online collaboration
$ws.Range("B7").Value = "Good"
$ws.Range("C7").Value = "Very active at dicord, Shared a lot of info about machine learning "

# Row 20 - Robert, Criterion 1 International Collaboration
$ws.Range("B20").Value = "Good"
$ws.Range("C20").Value = "Quick response, good at sharing info from lectors "

# Update the view scroll position / selection to match the target state
$win = $excel.ActiveWindow
$win.ScrollRow = 17
$win.ScrollColumn = 1
$ws.Range("C20").Select()
